$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab (SA-HW15.xpc -> SA) ---
$ws.Name = "SA"

# --- Tiny recomputed-precision tweaks in row 13 (last-ULP changes) ---
$ws.Range("D13").Value = 0.9973347867715998
$ws.Range("H13").Value = 0.9973347867715998
$ws.Range("K13").Value = 0.9955507390570452
$ws.Range("N13").Value = 0.9944665240692594

# --- Tiny recomputed-precision tweaks in row 15 (last-ULP changes) ---
$ws.Range("C15").Value = 0.9372133604321332
$ws.Range("G15").Value = 0.9372133604321332
$ws.Range("M15").Value = 0.9371583660735142
$ws.Range("O15").Value = 0.9986700295665724

# --- New row 16: results for the "HexGrid-60degTilt5degRes" scheme (index 14) ---
# Copy the formatting (bold / bordered / centered style) used by column A's
# index cells (e.g. A15) onto the new index cell A16 before filling values.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.10314550960926
$ws.Range("D16").Value = 0.8481667350877083
$ws.Range("E16").Value = 1.027015663102508
$ws.Range("F16").Value = 0.9440817023415943
$ws.Range("G16").Value = 1.10314550960926
$ws.Range("H16").Value = 0.8481667350877083
$ws.Range("I16").Value = 1.082398467621259
$ws.Range("J16").Value = 0.935262215100173
$ws.Range("K16").Value = 1.049690506732624
$ws.Range("L16").Value = 0.8793947169436527
$ws.Range("M16").Value = 1.10314550960926
$ws.Range("N16").Value = 0.937591199095108
$ws.Range("O16").Value = 0.9806024025352675
$ws.Range("P16").Value = 0.9836444395673473
